$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add I0 and IF headers in I1/J1, matching H1's style ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data columns I (I0) and J (IF) for rows 2-30 ---
$iValues = @(5,9,3,1,1,1,3,6,5,7,7,6,9,6,8,7,9,9,7,8,4,5,7,7,5,4,4,4,3)
$jValues = @(7,9,6,6,4,5,8,9,8,8,7,7,9,6,9,8,9,9,7,8,5,5,7,7,5,5,5,4,3)

for ($r = 0; $r -lt 29; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
